# HomeMove submit button added
#
# - Adds a new row to the "HomeMovePage" sheet with the key "submitButton"
#   and its xpath value, mirroring the existing Key/Value rows on that sheet.
# - Makes "HomeMovePage" the active sheet/tab (it was "HomePage" before),
#   with the selection left on A15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomeMovePage")

# Set the value cell (B) before the key cell (A) so that the new shared
# strings are appended to the shared-string table in the same order as
# the target workbook (xpath value first, then the "submitButton" key).
$ws.Range("B12").Value = "//android.widget.Button[@resource-id='com.Etisalat.ETIDA:id/btn_submit_request']"
$ws.Range("A12").Value = "submitButton"

# Match the formatting used by the other Value cells in column B.
$ws.Range("B12").Style = $ws.Range("B3").Style

# Make HomeMovePage the active sheet (this also clears tabSelected on the
# previously active sheet and updates the workbook's activeTab).
$ws.Activate()
$ws.Range("A15").Select()
